# Update TPM-derived values in the active worksheet (Gdf6-Bmpr1a).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 6.780879499999999
$ws.Range("N2").Value = 13.561759
$ws.Range("O2").Value = 0.09314755032665376
$ws.Range("P2").Value = 0.07273600820493056
$ws.Range("Q2").Value = 1.5748796065135
$ws.Range("R2").Value = 9.449277639080998
$ws.Range("S2").Value = 0.09314755032665376
$ws.Range("T2").Value = 0.07273600820493056

# Row 3
$ws.Range("O3").Value = 0.3322252662272683
$ws.Range("P3").Value = 0.389136476570504
$ws.Range("S3").Value = 0.3322252662272683
$ws.Range("T3").Value = 0.389136476570504

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.592065666666667
$ws.Range("N4").Value = 4.776197
$ws.Range("O4").Value = 0.02186987938794236
$ws.Range("P4").Value = 0.02561625701948875
$ws.Range("Q4").Value = 0.3697620272803332
$ws.Range("R4").Value = 3.327858245522999
$ws.Range("S4").Value = 0.02186987938794236
$ws.Range("T4").Value = 0.02561625701948875

# Row 5
$ws.Range("M5").Value = 25.158886
$ws.Range("N5").Value = 50.31777200000001
$ws.Range("O5").Value = 0.3456024546443489
$ws.Range("P5").Value = 0.2698701456828592
$ws.Range("Q5").Value = 5.843226750157999
$ws.Range("R5").Value = 35.059360500948
$ws.Range("S5").Value = 0.3456024546443489
$ws.Range("T5").Value = 0.2698701456828592

# Row 6
$ws.Range("M6").Value = 14.00046133333333
$ws.Range("N6").Value = 42.001384
$ws.Range("O6").Value = 0.1923214645892228
$ws.Range("P6").Value = 0.225266723235713
$ws.Range("Q6").Value = 3.251649146050666
$ws.Range("R6").Value = 29.264842314456
$ws.Range("S6").Value = 0.1923214645892228
$ws.Range("T6").Value = 0.225266723235713

# Row 7
$ws.Range("M7").Value = 1.079828666666667
$ws.Range("N7").Value = 3.239486
$ws.Range("O7").Value = 0.01483338482456395
$ws.Range("P7").Value = 0.01737438928650463
$ws.Range("Q7").Value = 0.2507934473193333
$ws.Range("R7").Value = 2.257141025874
$ws.Range("S7").Value = 0.01483338482456395
$ws.Range("T7").Value = 0.01737438928650463
